$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TIC_mas10_tablas")

# Footnote-marker style update on the "TIC_mas10_tablas" indicator labels:
# "*", "**", "***" -> "(i)", "(ii)", "(iii)", and the row-15 label is
# reworded from "Personal especialista TIC****" to a full sentence.
$ws.Range("A15").Value = "Porcentaje de personal especialista TIC sobre el total de personal"
$ws.Range("A6").Value  = "Empleados que teletrabajan regularmente(i)"
$ws.Range("A7").Value  = "Empresas con sitio/página web(ii)"
$ws.Range("A8").Value  = "Empresas que permiten teletrabajo(i)"
$ws.Range("A9").Value  = "Empresas que utilizan medios sociales(ii)"
$ws.Range("A10").Value = "Empresas que realizan analítica de datos internamente(ii)"
$ws.Range("A11").Value = "Empresas que compran servicios de cloud computing(ii)"
$ws.Range("A12").Value = "Empresas que emplean tecnologías de IA(ii)"
$ws.Range("A13").Value = "Empresas que emplean especialistas en TIC(i)"
$ws.Range("A14").Value = "Empresas con mujeres especialistas TIC(iii)"

# This sheet ("TIC_mas10_tablas") becomes the active tab/selection, taking
# over from "ID_estad_act_i+d_t1"; final cursor position is A15.
$ws.Activate() | Out-Null
$ws.Range("A15").Select() | Out-Null
